$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The source site re-scraped/re-sorted a handful of matches that share the
# same kickoff date (2023-11-25). Rows 134-137 (columns F:V -- everything
# except the row index / country / tournament / season / date columns,
# which stay put) end up holding a different match each, following a
# 4-cycle:
#   134 <- old136 (Progresul Spartac x Steaua Bucuresti)
#   135 <- old134 (Mioveni x Metaloglobus Bucharest)
#   136 <- old137 (CSM Resita x Alexandria)
#   137 <- old135 (Csikszereda M. Ciuc x CSM Slatina)
# A new match (row 141, Concordia x Unirea Slobozia) is also appended.
# ---------------------------------------------------------------------------

# Row 134: Mioveni x Metaloglobus Bucharest  ->  Progresul Spartac x Steaua Bucuresti
$ws.Range("F134").Value = "Progresul Spartac"
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = "Steaua Bucuresti"
$ws.Range("I134").Value = 6
$ws.Range("J134").Value = 7.1
$ws.Range("K134").Value = "23/11/2023 22:12"
$ws.Range("L134").Value = 11.17
$ws.Range("M134").Value = "25/11/2023 09:51"
$ws.Range("N134").Value = 4.56
$ws.Range("O134").Value = "23/11/2023 22:12"
$ws.Range("P134").Value = 5.61
$ws.Range("Q134").Value = "25/11/2023 09:51"
$ws.Range("R134").Value = 1.37
$ws.Range("S134").Value = "23/11/2023 22:12"
$ws.Range("T134").Value = 1.27
$ws.Range("U134").Value = "25/11/2023 09:51"
$ws.Range("V134").Value = "https://www.betexplorer.com/football/romania/liga-2/progresul-spartac-csa-steaua-bucuresti/0hpTJtCS/"

# Row 135: Csikszereda M. Ciuc x CSM Slatina  ->  Mioveni x Metaloglobus Bucharest
$ws.Range("F135").Value = "Mioveni"
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = "Metaloglobus Bucharest"
$ws.Range("I135").Value = 2
$ws.Range("J135").Value = 1.86
$ws.Range("K135").Value = "23/11/2023 22:12"
$ws.Range("L135").Value = 2.02
$ws.Range("M135").Value = "25/11/2023 09:51"
$ws.Range("N135").Value = 3.22
$ws.Range("O135").Value = "23/11/2023 22:12"
$ws.Range("P135").Value = 3.25
$ws.Range("Q135").Value = "25/11/2023 09:51"
$ws.Range("R135").Value = 3.93
$ws.Range("S135").Value = "23/11/2023 22:12"
$ws.Range("T135").Value = 3.91
$ws.Range("U135").Value = "25/11/2023 09:51"
$ws.Range("V135").Value = "https://www.betexplorer.com/football/romania/liga-2/mioveni-metaloglobus-bucharest/CWDoxORd/"

# Row 136: Progresul Spartac x Steaua Bucuresti  ->  CSM Resita x Alexandria
$ws.Range("F136").Value = "CSM Resita"
$ws.Range("G136").Value = 3
$ws.Range("H136").Value = "Alexandria"
$ws.Range("I136").Value = 1
$ws.Range("J136").Value = 1.68
$ws.Range("K136").Value = "23/11/2023 22:12"
$ws.Range("L136").Value = 1.63
$ws.Range("M136").Value = "25/11/2023 09:55"
$ws.Range("N136").Value = 3.61
$ws.Range("O136").Value = "23/11/2023 22:12"
$ws.Range("P136").Value = 3.87
$ws.Range("Q136").Value = "25/11/2023 09:58"
$ws.Range("R136").Value = 4.34
$ws.Range("S136").Value = "23/11/2023 22:12"
$ws.Range("T136").Value = 5.35
$ws.Range("U136").Value = "25/11/2023 09:57"
$ws.Range("V136").Value = "https://www.betexplorer.com/football/romania/liga-2/csm-resita-csm-alexandria/WzyaEKZk/"

# Row 137: CSM Resita x Alexandria  ->  Csikszereda M. Ciuc x CSM Slatina
$ws.Range("F137").Value = "Csikszereda M. Ciuc"
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = "CSM Slatina"
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 1.94
$ws.Range("K137").Value = "25/11/2023 00:42"
$ws.Range("L137").Value = 2.11
$ws.Range("M137").Value = "25/11/2023 09:58"
$ws.Range("N137").Value = 3.31
$ws.Range("O137").Value = "25/11/2023 00:42"
$ws.Range("P137").Value = 3.13
$ws.Range("Q137").Value = "25/11/2023 09:58"
$ws.Range("R137").Value = 4.04
$ws.Range("S137").Value = "25/11/2023 00:42"
$ws.Range("T137").Value = 3.81
$ws.Range("U137").Value = "25/11/2023 09:58"
$ws.Range("V137").Value = "https://www.betexplorer.com/football/romania/liga-2/miercurea-ciuc-csm-slatina/ITfOKMdM/"

# ---------------------------------------------------------------------------
# Append the new match as row 141. Duplicate the previous last row first so
# the new row inherits identical cell formatting/styles (bold+border index
# column, date-time number format, etc.), then overwrite with real values.
# ---------------------------------------------------------------------------
$ws.Range("A140:V140").Copy($ws.Range("A141"))

$ws.Range("A141").Value = 140
$ws.Range("B141").Value = "romania"
$ws.Range("C141").Value = "liga-2"
$ws.Range("D141").Value = "2023-2024"
$ws.Range("E141").Value = 45258.66666666666
$ws.Range("F141").Value = "Concordia"
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = "Unirea Slobozia"
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 1.79
$ws.Range("K141").Value = "27/11/2023 04:12"
$ws.Range("L141").Value = 1.95
$ws.Range("M141").Value = "28/11/2023 15:56"
$ws.Range("N141").Value = 3.18
$ws.Range("O141").Value = "27/11/2023 04:12"
$ws.Range("P141").Value = 2.92
$ws.Range("Q141").Value = "28/11/2023 15:56"
$ws.Range("R141").Value = 4.37
$ws.Range("S141").Value = "27/11/2023 04:12"
$ws.Range("T141").Value = 4.95
$ws.Range("U141").Value = "28/11/2023 15:56"
$ws.Range("V141").Value = "https://www.betexplorer.com/football/romania/liga-2/concordia-unirea-slobozia/WChCNOC3/"
